$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (tab name in workbook.xml)
$ws.Name = "Sheet1"

# Update header text in B1
$ws.Range("B1").Value = "Schwingungsdauer (s)"

# Build the bold / centered-top / thin-boxed header style on a single cell
# first (A1), then propagate it to B1 via a format-only copy so both cells
# end up sharing one cell style record instead of each accumulating their
# own chain of intermediate styles.
$a1 = $ws.Range("A1")
$a1.Font.Bold = $true
$a1.HorizontalAlignment = -4108  # xlCenter
$a1.VerticalAlignment = -4160    # xlTop
$a1.Borders.LineStyle = 1        # xlContinuous -> renders as a thin box border

$a1.Copy()
$b1 = $ws.Range("B1")
$b1.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update measurement values in column B
$ws.Range("B2").Value = 2.05
$ws.Range("B4").Value = 2.06
$ws.Range("B5").Value = 1.97
$ws.Range("B6").Value = 2.01
$ws.Range("B7").Value = 2
$ws.Range("B8").Value = 2.03
$ws.Range("B9").Value = 1.97
$ws.Range("B10").Value = 2.02
$ws.Range("B11").Value = 1.96
